$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.243.17"
$ws.Range("E2").Value = "  +3.47%  "
$ws.Range("D3").Value = "3.484.06"
$ws.Range("E3").Value = "  +3.66%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "409.16"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").Value = "131.98"
$ws.Range("E6").Value = "  +17.80%  "
$ws.Range("D7").Value = "3.478.87"
$ws.Range("E7").Value = "  +3.78%  "
$ws.Range("E8").Value = "  +3.09%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "0.700"
$ws.Range("E10").Value = "  +8.68%  "
$ws.Range("E11").Value = "  +31.57%  "
$ws.Range("D12").Value = "43.61"
$ws.Range("E12").Value = "  +9.85%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "4.029.99"
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "0.142"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("E15").Value = "  +4.63%  "
$ws.Range("D16").Value = "20.19"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").Value = "3.476.57"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "63.161.48"
$ws.Range("E18").Value = "  +3.82%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "0.0000144"
$ws.Range("E21").Value = "  +29.17%  "
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").Value = "82.69"
$ws.Range("E23").Value = "  +10.01%  "
$ws.Range("D24").Value = "13.16"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "313.00"
$ws.Range("E25").Value = "  +2.97%  "
$ws.Range("D26").Value = "3.16"
$ws.Range("E26").Value = "  -1.24%  "
$ws.Range("D27").Value = "30.59"
$ws.Range("E27").Value = "  +6.00%  "
$ws.Range("D28").Value = "8.19"
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "4.38"
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("D31").Value = "7.59"
$ws.Range("E31").Value = "  -2.66%  "
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").Value = "44.06"
$ws.Range("E33").Value = "  +11.56%  "
$ws.Range("D34").Value = "11.89"
$ws.Range("E34").Value = "  +3.80%  "
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("D38").Value = "52.68"
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").Value = "3.59"
$ws.Range("E39").Value = "  +5.29%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("E43").Value = "  +4.07%  "
$ws.Range("D44").Value = "136.97"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "17.57"
$ws.Range("E45").Value = "  +3.94%  "
$ws.Range("E46").Value = "  -4.68%  "
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").Value = "22.31"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").Value = "3.825.56"
$ws.Range("E50").Value = "  +3.59%  "
$ws.Range("D51").Value = "2.187.61"
$ws.Range("E51").Value = "  -0.02%  "
